$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new 2022 column (M), copying the formatting of the matching
# existing cells before filling in the new values.
$ws.Range("K4").Copy()
$ws.Range("M4").PasteSpecial(-4122)

$ws.Range("L5").Copy()
$ws.Range("M5").PasteSpecial(-4122)

$ws.Range("L6").Copy()
$ws.Range("M6").PasteSpecial(-4122)

$ws.Range("M4").Value = 2022
$ws.Range("M5").Value = 2.2
$ws.Range("M6").Value = 1.2

# Update the selection to match the new active cell
$ws.Range("M10").Select()
